$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original styles for D2:D51 and E2:E51 so we can restore them
# after forcing a Text number format (needed so values like "1.002" or
# "0.2190" are stored as literal text instead of being auto-converted to
# numbers by Excel's smart-parsing, which would drop trailing zeros /
# collapse "30.358.05" style multi-dot strings).
$rngD = $ws.Range("D2:D51")
$rngE = $ws.Range("E2:E51")
$styleD = $rngD.Style
$styleE = $rngE.Style
$rngD.NumberFormat = "@"
$rngE.NumberFormat = "@"

$ws.Range("D2").Value = "30.358.05"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "2.013.43"
$ws.Range("E3").Value = "  +5.49%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "325.01"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.5154"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("D8").Value = "0.4273"
$ws.Range("E8").Value = "  +5.48%  "
$ws.Range("D9").Value = "0.08701"
$ws.Range("E9").Value = "  +5.23%  "
$ws.Range("D10").Value = "43.25"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("D11").Value = "1.133"
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("D12").Value = "24.76"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").Value = "2.009.30"
$ws.Range("E13").Value = "  +5.46%  "
$ws.Range("D14").Value = "6.581"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "7.460"
$ws.Range("E15").Value = "  +3.49%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "94.54"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "0.00001112"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "18.89"
$ws.Range("E20").Value = "  +4.19%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "6.204"
$ws.Range("D23").Value = "30.419.95"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "11.81"
$ws.Range("E24").Value = "  +4.45%  "
$ws.Range("D25").Value = "2.242"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").Value = "2.245.26"
$ws.Range("E26").Value = "  +5.67%  "
$ws.Range("D27").Value = "22.37"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "162.81"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "2.415"
$ws.Range("E29").Value = "  +5.20%  "
$ws.Range("D30").Value = "131.37"
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("D31").Value = "1.137"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").Value = "6.068"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").Value = "3.831"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").Value = "1.363"
$ws.Range("E35").Value = "  +14.06%  "
$ws.Range("D36").Value = "0.02529"
$ws.Range("E36").Value = "  +3.72%  "
$ws.Range("D37").Value = "5.469"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").Value = "0.06642"
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("D39").Value = "12.33"
$ws.Range("E39").Value = "  +8.39%  "
$ws.Range("D40").Value = "9.053"
$ws.Range("E40").Value = "  +4.10%  "
$ws.Range("D41").Value = "0.2190"
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("D42").Value = "0.6651"
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "13.60"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").Value = "0.6162"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "2.184"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").Value = "3.660"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").Value = "1.264"
$ws.Range("E49").Value = "  +4.93%  "
$ws.Range("D50").Value = "124.15"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").Value = "80.83"
$ws.Range("E51").Value = "  +2.68%  "

# Restore original (unstyled) formatting
$rngD.Style = $styleD
$rngE.Style = $styleE
